$d = $word.ActiveDocument

# Color used for highlighted metrics: 2C3E50 (RGB) -> stored as BGR-order
# integer for Word's Font.Color (wdColor) property.
$metricColor = 5258796  # 0x2C3E50 in RGB becomes 0x503E2C in BGR == 5258796

function Apply-MetricHighlight($ParaIndex, $Metrics) {
    $p = $d.Paragraphs($ParaIndex)
    $paraStart = $p.Range.Start
    $cursor = 0

    foreach ($metric in $Metrics) {
        $t = $p.Range.Text
        $idx = $t.IndexOf($metric, $cursor)
        if ($idx -lt 0) {
            continue
        }
        $rangeStart = $paraStart + $idx
        $rangeEnd = $rangeStart + $metric.Length
        $r = $d.Range($rangeStart, $rangeEnd)
        $r.Font.Bold = 1
        $r.Font.Color = $metricColor
        $cursor = $idx + $metric.Length
    }
}

# "Discovered systematic race coding errors ... from 23% to 64%"
Apply-MetricHighlight 10 @("23%", "64%")

# "Utilized advanced sampling methods ... from ±4.2% to ±2.1%, ... from 71% to 87% ..."
Apply-MetricHighlight 12 @([char]0x00B1 + "4.2%", [char]0x00B1 + "2.1%", "71%", "87%")

# "Trigonometric algorithm ... by 73.5%, saving ... $4.7M and enabling ..."
Apply-MetricHighlight 13 @("73.5%", "`$4.7M")

# "Built real-time FEC analysis systems ... valued over $2 trillion"
Apply-MetricHighlight 14 @("`$2")

# "Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%"
Apply-MetricHighlight 50 @("73.5%")

# "$4.7M savings enabled nonprofit access"
Apply-MetricHighlight 51 @("`$4.7M")

# "178% accuracy improvement in racial classification algorithms"
Apply-MetricHighlight 53 @("178%")
